{"js": "// Update the SUMMARY paragraph of the resume:\n//  - \" R&D experience.\"            -> \" research background.\"\n//  - \" with a math and physics focus.\" -> \" with a modeling and analysis focus.\"\n//  - \"Versed in missile design, GNC, radar, stereo vision, airfoil\"\n//                                   -> \"Versed in controls, stereo vision, missile design, data analysis, airfoil\"\nconst body = context.document.body;\n\nconst replacements = [\n  [\" R&D experience.\", \" research background.\"],\n  [\" with a math and physics focus.\", \" with a modeling and analysis focus.\"],\n  [\n    \"Versed in missile design, GNC, radar, stereo vision, airfoil\",\n    \"Versed in controls, stereo vision, missile design, data analysis, airfoil\",\n  ],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + find);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the SUMMARY paragraph of the resume:\n#  - \" R&D experience.\"            -> \" research background.\"\n#  - \" with a math and physics focus.\" -> \" with a modeling and analysis focus.\"\n#  - \"Versed in missile design, GNC, radar, stereo vision, airfoil\"\n#                                   -> \"Versed in controls, stereo vision, missile design, data analysis, airfoil\"\n\n$d = $word.ActiveDocument\n\n# Applied right-to-left (latest match first) so each Find/Replace only\n# re-merges the runs that follow it, leaving earlier, untouched text runs\n# (and their formatting) intact.\n$replacements = @(\n    @(\"Versed in missile design, GNC, radar, stereo vision, airfoil\", \"Versed in controls, stereo vision, missile design, data analysis, airfoil\"),\n    @(\" with a math and physics focus.\", \" with a modeling and analysis focus.\"),\n    @(\" R&D experience.\", \" research background.\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $findText\"\n    }\n}\n"}
